$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.482.99"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = "'  +0.46%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Range("D3").Value = "'1.797.45"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = "'  -0.19%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Range("D5").Value = "'316.87"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "'  +0.32%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Range("D7").Value = "'0.5415"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "'  -1.62%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Range("D8").Value = "'0.3782"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "'  -1.22%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Range("D9").Value = "'0.07496"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "'  -0.80%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Range("D10").Value = "'41.89"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "'  -2.22%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Range("D11").Value = "'1.107"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "'  -1.31%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Range("E12").Value = "'  -0.17%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Range("D13").Value = "'20.67"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "'  -2.10%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Range("D14").Value = "'6.156"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "'  -0.40%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Range("D15").Value = "'7.297"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "'  +0.07%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Range("D16").Value = "'1.796.39"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "'  -0.26%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Range("D17").Value = "'89.61"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "'  -1.28%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Range("D18").Value = "'0.00001066"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "'  +0.04%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Range("D19").Value = "'0.06503"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "'  +0.73%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Range("E20").Value = "'  +1.49%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Range("E21").Value = "'  +0.00%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Range("D22").Value = "'5.944"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "'  -0.36%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Range("D23").Value = "'28.486.75"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = "'  +0.33%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Range("D24").Value = "'11.10"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "'  -1.01%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Range("E25").Value = "'  -1.26%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Range("D26").Value = "'159.33"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E27").Value = "'  -1.06%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Range("D28").Value = "'2.000.74"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "'  -0.58%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Range("E29").Value = "'  -4.18%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Range("D30").Value = "'122.76"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "'  -0.91%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Range("D31").Value = "'1.113"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "'  -4.07%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Range("D32").Value = "'0.1053"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "'  +2.98%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Range("E33").Value = "'  -1.60%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Range("D34").Value = "'3.650"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "'  -0.65%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Range("D35").Value = "'0.2269"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "'  -0.53%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Range("D36").Value = "'0.06459"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "'  +3.50%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Range("D37").Value = "'0.02295"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "'  -0.68%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Range("D38").Value = "'8.610"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "'  -2.90%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Range("D39").Value = "'5.029"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "'  +0.61%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Range("D40").Value = "'0.6201"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "'  -2.30%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Range("E41").Value = "'  -3.04%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Range("D42").Value = "'1.451"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "'  +4.78%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Range("D43").Value = "'1.192"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "'  +2.33%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Range("E44").Value = "'  -0.03%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Range("D45").Value = "'13.36"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "'  -0.61%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Range("D46").Value = "'3.686"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "'  +0.48%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Range("D47").Value = "'0.5824"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "'  -2.10%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Range("D48").Value = "'127.28"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "'  +3.19%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Range("D49").Value = "'1.207"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "'  +5.55%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Range("D50").Value = "'1.949"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "'  -0.27%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Range("D51").Value = "'0.06892"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "'  -0.14%  "
$ws.Cells.Item(51, 5).Style = "Normal"
